# Update weekly Fruta/Hortaliza (Guayaba) records: reassign Fecha (D), Volumen (M)
# and the associated price columns (N, O, P, S) per row, matching the new data
# scrape for the period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value  = 44424
$ws.Cells.Item(2, 13).Value = 50
$ws.Cells.Item(2, 14).Value = 1200
$ws.Cells.Item(2, 15).Value = 1200
$ws.Cells.Item(2, 16).Value = 1200
$ws.Cells.Item(2, 19).Value = 1200

# Row 3
$ws.Cells.Item(3, 4).Value  = 44343
$ws.Cells.Item(3, 14).Value = 1300
$ws.Cells.Item(3, 15).Value = 1300
$ws.Cells.Item(3, 16).Value = 1300
$ws.Cells.Item(3, 19).Value = 1300

# Row 4
$ws.Cells.Item(4, 4).Value  = 44476
$ws.Cells.Item(4, 13).Value = 80

# Row 5
$ws.Cells.Item(5, 4).Value  = 44418
$ws.Cells.Item(5, 13).Value = 40
$ws.Cells.Item(5, 14).Value = 1200
$ws.Cells.Item(5, 15).Value = 1200
$ws.Cells.Item(5, 16).Value = 1200
$ws.Cells.Item(5, 19).Value = 1200

# Row 6
$ws.Cells.Item(6, 4).Value  = 44473
$ws.Cells.Item(6, 13).Value = 120
$ws.Cells.Item(6, 14).Value = 1200
$ws.Cells.Item(6, 15).Value = 1200
$ws.Cells.Item(6, 16).Value = 1200
$ws.Cells.Item(6, 19).Value = 1200

# Row 7
$ws.Cells.Item(7, 4).Value  = 44438
$ws.Cells.Item(7, 13).Value = 60

# Row 9
$ws.Cells.Item(9, 4).Value  = 44435
$ws.Cells.Item(9, 13).Value = 130
$ws.Cells.Item(9, 14).Value = 1300
$ws.Cells.Item(9, 15).Value = 1300
$ws.Cells.Item(9, 16).Value = 1300
$ws.Cells.Item(9, 19).Value = 1300

# Row 11
$ws.Cells.Item(11, 4).Value = 44417

# Row 12
$ws.Cells.Item(12, 4).Value  = 44431
$ws.Cells.Item(12, 13).Value = 100

# Row 13
$ws.Cells.Item(13, 4).Value  = 44432
$ws.Cells.Item(13, 13).Value = 30
$ws.Cells.Item(13, 14).Value = 1300
$ws.Cells.Item(13, 15).Value = 1300
$ws.Cells.Item(13, 16).Value = 1300
$ws.Cells.Item(13, 19).Value = 1300
